# 07-07-22 - Se actualizan cifras del 6 y 7 de Julio 2022
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Julio")
$ws.Activate()

# Update selection to C4
$ws.Range("C4").Select()

# Row 7 (6 de julio) values: B..J
$ws.Cells.Item(7, 2).Value = 78
$ws.Cells.Item(7, 3).Value = 0
$ws.Cells.Item(7, 4).Value = 3
$ws.Cells.Item(7, 5).Value = 4
$ws.Cells.Item(7, 6).Value = 0
$ws.Cells.Item(7, 7).Value = 1
$ws.Cells.Item(7, 8).Value = 0
$ws.Cells.Item(7, 9).Value = 0
$ws.Cells.Item(7, 10).Value = 0

# Row 8 (7 de julio) values: B..J
$ws.Cells.Item(8, 2).Value = 92
$ws.Cells.Item(8, 3).Value = 2
$ws.Cells.Item(8, 4).Value = 2
$ws.Cells.Item(8, 5).Value = 0
$ws.Cells.Item(8, 6).Value = 0
$ws.Cells.Item(8, 7).Value = 0
$ws.Cells.Item(8, 8).Value = 0
$ws.Cells.Item(8, 9).Value = 0
$ws.Cells.Item(8, 10).Value = 0

# Row 24: F24 gets an underline style applied (no value)
$ws.Cells.Item(24, 6).Font.Underline = $true
